$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "success" column (D) indicating whether a row's list is the
# most frequent one (count = 377 -> success = 1, everything else -> 0).

# Header cell D1, formatted like the other header cells (B1/C1).
$ws.Range("D1").Value = "success"
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# D2 is the winning/most frequent list -> success = "1"
$ws.Range("D2").Value = "'1"

# D3:D17 are the rest -> success = "0"
for ($r = 3; $r -le 17; $r++) {
    $ws.Cells.Item($r, 4).Value = "'0"
}

# Values above are entered with a leading apostrophe so Excel keeps them
# as text (matching the original "0"/"1" shared strings) instead of
# converting them to numbers. Re-apply the plain data-cell formatting
# (same as column B) so the quote-prefix flag doesn't linger as a style.
$ws.Range("B2").Copy()
$ws.Range("D2:D17").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
